$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76's date cell currently uses the "last row" date-only format (style 3).
# Since we're appending a new last row (77), row 76 reverts to the normal
# datetime format used by all the other data rows (style 2, same as A75).
$ws.Range("A76").NumberFormat = $ws.Range("A75").NumberFormat

# Add the new row of data (row 77), with the date-only format moved here.
$ws.Range("A77").Value = 45817
$ws.Range("B77").Value = 328
$ws.Range("C77").Value = 327
$ws.Range("D77").Value = 332

$ws.Range("A77").NumberFormat = "YYYY-MM-DD"
